$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 3 new daily rows (07, 08 and 12 October 2021), matching the same
# monetary policy / standing facility rates as the prior row (193).
$newRows = @(
    @{ Row = 194; Date = "07-10-2021"; TPM = 1.5; FPL = 1.75; FPD = 1.25 },
    @{ Row = 195; Date = "08-10-2021"; TPM = 1.5; FPL = 1.75; FPD = 1.25 },
    @{ Row = 196; Date = "12-10-2021"; TPM = 1.5; FPL = 1.75; FPD = 1.25 }
)

# Temporarily format the new column-A cells as text so the "dd-mm-yyyy"
# looking strings are kept as literal text (like the rest of column A)
# instead of being auto-converted into date serial numbers, then drop the
# formatting again so the cells end up unstyled - matching every other
# date cell already in the sheet.
$dateRange = $ws.Range("A194:A196")
$dateRange.NumberFormat = "@"

foreach ($row in $newRows) {
    $ws.Cells.Item($row.Row, 1).Value = $row.Date
    $ws.Cells.Item($row.Row, 2).Value = $row.TPM
    $ws.Cells.Item($row.Row, 3).Value = $row.FPL
    $ws.Cells.Item($row.Row, 4).Value = $row.FPD
}

$dateRange.ClearFormats()
